$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201, pushing existing rows 201:310 down to 202:311.
$ws.Rows("201:201").Insert()

# Populate the newly inserted row 201 with the new data record.
$ws.Range("A201").Value = 9
$ws.Range("B201").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C201").Value = "Metropolitana"
$ws.Range("D201").Value = 44917
$ws.Range("E201").Value = 13
$ws.Range("F201").Value = 100112026
$ws.Range("G201").Value = "Haba"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 90
$ws.Range("K201").Value = 16000
$ws.Range("L201").Value = 18000
$ws.Range("M201").Value = 17000
$ws.Range("N201").Value = "$/saco 25 kilos"
$ws.Range("O201").Value = "Carahue"
$ws.Range("P201").Value = 680
$ws.Range("Q201").Value = 25
$ws.Range("R201").Value = "Hortaliza"
